# Applies the LOM3272.docx restructuring described in the commit diff.
# The edit reshuffles paragraph *contents* (text/runs/paragraph style) across
# the document while the number of paragraphs (20) is unchanged, so each
# affected paragraph is rewritten in place via Range.InsertXML with the exact
# OOXML the target revision contains for that paragraph.

$d = $word.ActiveDocument

$expectedCount = 20
if ($d.Paragraphs.Count -ne $expectedCount) {
    throw "Expected $expectedCount paragraphs, found $($d.Paragraphs.Count); aborting to avoid corrupting the document."
}

# Para 6 ('Objetivos' body, PT): now holds the 'Programa resumido' PT paragraph text
$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Os cursos de engenharia física, respectivos projetos pedagógicos e seus componentes curriculares, incluindo TCC, estágio obrigatório, Projetos de Extensão Curricularizados, Atividades Acadêmicas Complementares e Atividades extracurriculares. Identificação e aderência do estudante com o curso e com a profissão escolhida. O curso superior, a transição adolescente/jovem adulto e os desafios nos projetos de vida do estudante no início da graduação. Relação entre as disciplinas e o conhecimento a ser aplicado. Competências e habilidades desenvolvidas no seu curso de engenharia. Dimensões acadêmicas, socioculturais e científicas. Diversidade e inclusão. Organização dos estudos.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(6).Range.InsertXML($xml6)

# Para 7 (italic EN under 'Objetivos'): now holds the 'Programa resumido' EN paragraph text
$xml7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>Engineering physics courses, respective pedagogical projects and their curricular components, including TCC, mandatory internship, Curricular Extension Projects, Complementary Academic Activities and Extracurricular Activities. Identification and adherence of the student with the course and with the chosen profession. The college course, the adolescent/young adult transition and the challenges in the student's life projects at the beginning of the undergraduate program. Relationship between the disciplines and the knowledge to be applied. Competencies and skills developed in your engineering course. Academic, sociocultural and scientific dimensions. Diversity and inclusion. Organization of studies.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(7).Range.InsertXML($xml7)

# Para 9 (List Bullet under 'Docente(s) Responsavel(eis)'): now holds the old PT objectives
# text plus the old 'Programa' PT text, joined by a manual line break
$xml9 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Orientar os estudantes no início de sua trajetória universitária no curso de graduação em Engenharia XX na EEL-USP de modo que o estudante seja capaz de a) identificar as oportunidades acadêmicas e as particularidades do seu curso; b) reconhecer, sob acompanhamento de um tutor, eventuais dificuldades ao longo do curso e compreender mecanismos para que estas sejam superadas, conduzindo o curso com o sucesso desejado; c) desenvolver habilidades técnicas e emocionais, ampliando as perspectivas de formação profissional por meio de atividades e encontros sistematizados.</w:t><w:br/></w:r><w:r><w:t>Apresentação dos programas e serviços oferecidos pela USP voltados aos estudantes e das oportunidades de realizar trabalhos extracurriculares. A dinâmica das aulas, ferramentas de interação. Desenvolvimento de atividades de grupo, com objetivo de desenvolver habilidades sócio-comportamentais através de colaboração em temas do curso relacionados à profissão escolhida. Áreas de atuação do curso de engenharia, competências e habilidades a serem desenvolvidas. Interdisciplinaridade e a relação entre as disciplinas e o conhecimento a ser aplicado. Planejamento de estudos. Formas de estudar e aprender.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(9).Range.InsertXML($xml9)

# Para 11 (body under 'Programa resumido'): now holds the old 'Metodo' text
$xml11 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Atividades realizadas na forma de dinâmicas de grupos, utilização de vídeos, textos, roda de discussão e/ou elaboração de painéis. Participação em encontros de orientação promovidos pelo Programa de Tutoria Acadêmica e a realização de atividades propostas pelo tutor/monitor/mentor, incluindo trabalhos em equipe e estudos dirigidos.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(11).Range.InsertXML($xml11)

# Para 12 (italic, under 'Programa resumido'): now holds the EN objectives paragraph
$xml12 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>To guide students at the beginning of their university career in the undergraduate course in Engineering XX at EEL-USP so that the student is able to: a) identify the academic opportunities and particularities of their course; b) recognize, under the supervision of a tutor, any difficulties throughout the course and understand mechanisms for them to be overcome, conducting the course with the desired success; c) develop technical and emotional skills, broadening the perspectives of professional training through systematized activities and meetings.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(12).Range.InsertXML($xml12)

# Para 14 (body under 'Programa'): now holds the old 'Criterio' text
$xml14 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Participação ativa nos encontros, apresentação de estudos/pesquisa e de trabalhos realizados durante a disciplina, colaboração e engajamento nas atividades da disciplina. O estudante deverá entregar um relatório final para a disciplina. A nota final é dada pela média ponderada das notas obtidas nas diversas atividades propostas.</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(14).Range.InsertXML($xml14)

# Para 16 (List Bullet under 'Avaliacao'): Metodo/Criterio/Norma runs now hold the
# 'Nao se aplica..' text, the bibliography entries, and the first docente name
$xml16 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Método: </w:t></w:r><w:r><w:t>Não se aplica..</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Critério: </w:t></w:r><w:r><w:t>[1] Peddy, S. The art of mentoring – Lead, follow and get out of the way. Houston: Bullion Books, 2001.</w:t><w:br/><w:t>[2] Zachary, L. J. The Mentor’s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promoção do bem-estar em estudantes do ensino superior. In: Programa de Monitorização e Tutorado: oito anos a promover a integração e o sucesso académico no IST. Lisboa: IST Press, 2011. p. 19-27.</w:t><w:br/><w:t>[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.</w:t><w:br/><w:t>[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.</w:t><w:br/><w:t>[5] Diretrizes Curriculares Nacionais para os cursos de graduação em Engenharia. Ministério da Educação. CNE/CES, 2019.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Norma de recuperação: </w:t></w:r><w:r><w:t>5817692 - Katia Cristiane Gandolpho Candioto</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(16).Range.InsertXML($xml16)

# Para 18 (body under 'Bibliografia'): now holds the second docente name
$xml18 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>1176388 - Luiz Tadeu Fernandes Eleno</w:t></w:r></w:p>
'@
[void]$d.Paragraphs.Item(18).Range.InsertXML($xml18)

Write-Host "LOM3272.docx restructuring applied."